$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: AD1, AE1, AF1 - copy formatting from an existing header cell, then set values
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Data rows 2-49: AD=67 (Wins), AE=95 (Losses), AF=0 (Ties)
for ($r = 2; $r -le 49; $r++) {
    $ws.Cells.Item($r, 30).Value = 67
    $ws.Cells.Item($r, 31).Value = 95
    $ws.Cells.Item($r, 32).Value = 0
}
